$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "개혁과 안티와 대응"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/revolution-anti-handling/#utm_source=rss&utm_medium=rss&utm_campaign=revolution-anti-handling"

$ws.Range("D52").Value = "[R] 3판 맛보기) 유니코드 문자열"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2560&utm_source=rss&utm_medium=rss&utm_campaign=r-3%25ed%258c%2590-%25eb%25a7%259b%25eb%25b3%25b4%25ea%25b8%25b0-%25ec%259c%25a0%25eb%258b%2588%25ec%25bd%2594%25eb%2593%259c-%25eb%25ac%25b8%25ec%259e%2590%25ec%2597%25b4"
